$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the two newly-reserved items in column E with "Y"
$ws.Range("E15").Value = "Y"
$ws.Range("E16").Value = "Y"

# Move the active selection/cursor to H22 (matches the saved cursor position in the file)
$ws.Activate()
$ws.Range("H22").Select()
